$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'27.422.67"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  +3.85%  "
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "'1.838.12"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  +3.50%  "
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(4, 4).Value = "'1.027"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  +2.56%  "
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "'317.93"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  +4.18%  "
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  +2.29%  "
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 4).Value = "'0.4361"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  +2.56%  "
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(8, 4).Value = "'0.3724"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'  +2.88%  "
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(9, 4).Value = "'0.07362"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  +2.62%  "
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 4).Value = "'0.8730"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  +3.92%  "
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 4).Value = "'21.31"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  +4.17%  "
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 4).Value = "'1.891.26"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  +5.42%  "
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'5.463"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  +4.09%  "
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 4).Value = "'6.689"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  +3.65%  "
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 4).Value = "'0.07111"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  +3.87%  "
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 4).Value = "'82.33"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  +4.33%  "
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 4).Value = "'1.030"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  +2.40%  "
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 4).Value = "'0.000008994"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  +3.70%  "
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(19, 4).Value = "'1.024"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  +2.05%  "
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 4).Value = "'15.39"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  +2.81%  "
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 4).Value = "'27.433.30"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  +3.85%  "
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 4).Value = "'5.232"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  +2.92%  "
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 4).Value = "'11.15"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  +0.81%  "
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 4).Value = "'156.86"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  +2.95%  "
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'1.892"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  +4.63%  "
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 4).Value = "'18.57"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  +2.96%  "
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 4).Value = "'5.231"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  +3.01%  "
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 4).Value = "'1.917"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  +5.87%  "
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 4).Value = "'115.77"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  +1.48%  "
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 4).Value = "'0.09040"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  +1.96%  "
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 4).Value = "'1.200"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  +6.90%  "
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 4).Value = "'0.7581"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  +4.49%  "
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 4).Value = "'4.467"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  +3.23%  "
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 4).Value = "'2.868"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  +4.51%  "
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 4).Value = "'1.026"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  +2.48%  "
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 4).Value = "'1.147"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  +4.41%  "
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 4).Value = "'0.01961"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  +4.11%  "
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 4).Value = "'0.05246"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  +2.01%  "
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(39, 4).Value = "'0.5154"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  +4.77%  "
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(40, 4).Value = "'2.786"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  +7.28%  "
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 4).Value = "'0.1660"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  +2.80%  "
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 4).Value = "'6.559"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  +2.78%  "
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(43, 4).Value = "'8.495"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  +5.92%  "
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(44, 4).Value = "'108.59"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  +3.55%  "
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 4).Value = "'10.56"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  +3.45%  "
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 4).Value = "'1.027"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  +2.52%  "
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 4).Value = "'1.678"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  +2.49%  "
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 4).Value = "'0.4621"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  +3.53%  "
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 2).Value = "'RenderToken"
$ws.Cells.Item(49, 2).Style = "Normal"
$ws.Cells.Item(49, 3).Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(49, 3).Style = "Normal"
$ws.Cells.Item(49, 4).Value = "'1.891"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  +9.55%  "
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 2).Value = "'Cronos"
$ws.Cells.Item(50, 2).Style = "Normal"
$ws.Cells.Item(50, 3).Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 3).Style = "Normal"
$ws.Cells.Item(50, 4).Value = "'0.06307"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  +2.08%  "
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 4).Value = "'39.42"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  +7.08%  "
$ws.Cells.Item(51, 5).Style = "Normal"
